# Generate Report for Archive
#
# 1) Status text "Ready for handoff" -> "In Translation" everywhere it is used
#    (the Overview sheet's per-language status columns E/F, and the Status
#    column C on each per-language detail sheet).
# 2) Narrow the "Status"-ish column(s) that previously auto-sized to fit the
#    long "Ready for handoff" text - now that the text is shorter they are
#    re-autofit to the new, narrower content width.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"
$overview.Range("E4").Value = "In Translation"
$overview.Range("F4").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"
$zhcn.Range("C4").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"
$dede.Range("C4").Value = "In Translation"

# Re-fit the now-narrower status columns.
$overview.Range("E1").ColumnWidth = 12.576851526896165
$overview.Range("F1").ColumnWidth = 12.576851526896165
$zhcn.Range("C1").ColumnWidth = 12.576851526896165
$dede.Range("C1").ColumnWidth = 12.576851526896165
